$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.083.23'
$ws.Range("E2").Value = '  -1.22%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.318.71'
$ws.Range("E3").Value = '  +1.55%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '186.78'
$ws.Range("E5").Value = '  +1.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '577.07'
$ws.Range("E6").Value = '  -0.57%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.604'
$ws.Range("E8").Value = '  +0.22%  '
$ws.Range("E9").Value = '  -0.42%  '
$ws.Range("E10").Value = '  +1.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.413'
$ws.Range("E11").Value = '  +1.14%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.892.97'
$ws.Range("E12").Value = '  +1.44%  '
$ws.Range("E13").Value = '  -0.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.45'
$ws.Range("E14").Value = '  +0.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '67.329.20'
$ws.Range("E15").Value = '  -0.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000167'
$ws.Range("E16").Value = '  -0.65%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.321.33'
$ws.Range("E17").Value = '  +1.49%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '445.07'
$ws.Range("E18").Value = '  +10.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.71'
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.54'
$ws.Range("E20").Value = '  +1.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.79'
$ws.Range("E21").Value = '  +3.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.47'
$ws.Range("E22").Value = '  +3.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.514'
$ws.Range("E24").Value = '  +1.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.442.21'
$ws.Range("E25").Value = '  +1.02%  '
$ws.Range("E26").Value = '  +1.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.188'
$ws.Range("E27").Value = '  +0.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.11'
$ws.Range("E28").Value = '  -3.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("E30").Value = '  +1.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.91'
$ws.Range("E31").Value = '  +1.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.34'
$ws.Range("E32").Value = '  -2.22%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.24'
$ws.Range("E34").Value = '  -0.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.80'
$ws.Range("E35").Value = '  -1.28%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.53'
$ws.Range("E36").Value = '  +5.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '162.72'
$ws.Range("E37").Value = '  -0.98%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '27.59'
$ws.Range("E38").Value = '  +1.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.85'
$ws.Range("E39").Value = '  -2.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.790'
$ws.Range("E40").Value = '  -1.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.47'
$ws.Range("E41").Value = '  -0.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.743.03'
$ws.Range("E42").Value = '  +2.46%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.31'
$ws.Range("E43").Value = '  -0.41%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '24.89'
$ws.Range("E44").Value = '  +1.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0673'
$ws.Range("E45").Value = '  -0.67%  '
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.17'
$ws.Range("E46").Value = '  -1.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.41'
$ws.Range("E47").Value = '  -0.65%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '329.76'
$ws.Range("E48").Value = '  -1.55%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0275'
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.994'
$ws.Range("E50").Value = '  +2.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.27'
$ws.Range("E51").Value = '  -0.41%  '
